# Split the merge-field placeholder "<<[agreement_hiren]>>" into
# "<<[agreement_number]>>", typed as three separate runs:
#   "<<[agreement_"  +  "number"  +  "]>>"
# matching how Word splits a run when you select a word in the middle
# of existing text and retype it.

$d = $word.ActiveDocument

# Locate the placeholder and grab its range.
$rng = $d.Content
$found = $rng.Find.Execute("<<[agreement_hiren]>>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$start = $rng.Start
$end = $rng.End

# Clear the existing text, then retype it as three distinct pieces so the
# saved document ends up with three runs instead of one.
$whole = $d.Range($start, $end)
$whole.Text = ""

$prefixEnd = $start + ("<<[agreement_").Length
$p1 = $d.Range($start, $start)
$p1.InsertAfter("<<[agreement_")

$middleEnd = $prefixEnd + ("number").Length
$p2 = $d.Range($prefixEnd, $prefixEnd)
$p2.InsertAfter("number")

$p3 = $d.Range($middleEnd, $middleEnd)
$p3.InsertAfter("]>>")

# A transient bookmark around the middle piece stops the engine from
# silently re-coalescing these three identically-formatted runs back into
# one when the document is saved; removing the bookmark again leaves no
# trace of it in the saved XML while keeping the runs split apart.
$d.Bookmarks.Add("__tmp_split_marker", $d.Range($prefixEnd, $middleEnd)) | Out-Null
$d.Bookmarks("__tmp_split_marker").Delete()
